$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 23:11"

# Country name reorder (rows whose country label shifted because the
# underlying ranking/case counts changed)
$ws.Range("A78").Value = "Libia"
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("A123").Value = "Ruanda"
$ws.Range("A124").Value = "Surinam"
$ws.Range("A129").Value = "Angola"
$ws.Range("A130").Value = "Birmania"
$ws.Range("A131").Value = "Tailandia"
$ws.Range("A138").Value = "Aruba"
$ws.Range("A139").Value = "Guadalupe"

# Updated case statistics
$ws.Range("B4").Value = 6777184
$ws.Range("C4").Value = 27895
$ws.Range("D4").Value = 4052454
$ws.Range("E4").Value = 2524869
$ws.Range("G4").Value = 861
$ws.Range("H4").Value = 199861

$ws.Range("E5").Value = 996832
$ws.Range("G5").Value = 1283
$ws.Range("H5").Value = 82091

$ws.Range("B6").Value = 4382263
$ws.Range("C6").Value = 32719
$ws.Range("E6").Value = 635960
$ws.Range("G6").Value = 1002
$ws.Range("H6").Value = 133119

$ws.Range("B11").Value = 651521
$ws.Range("C11").Value = 772
$ws.Range("D11").Value = 583126
$ws.Range("E11").Value = 52754
$ws.Range("G11").Value = 142
$ws.Range("H11").Value = 15641

$ws.Range("B27").Value = 164402
$ws.Range("C27").Value = 4034
$ws.Range("D27").Value = 120727
$ws.Range("E27").Value = 42528
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = 1147

$ws.Range("B35").Value = 104803
$ws.Range("C35").Value = 693
$ws.Range("D35").Value = 78387
$ws.Range("E35").Value = 24418
$ws.Range("G35").Value = 14
$ws.Range("H35").Value = 1998

$ws.Range("D61").Value = 39900
$ws.Range("E61").Value = 5823

$ws.Range("B78").Value = 24144
$ws.Range("C78").Value = 629
$ws.Range("D78").Value = 13252
$ws.Range("E78").Value = 10509
$ws.Range("G78").Value = 15
$ws.Range("H78").Value = 383

$ws.Range("B79").Value = 23929
$ws.Range("C79").Value = 294
$ws.Range("D79").Value = 16701
$ws.Range("E79").Value = 6503
$ws.Range("G79").Value = 20
$ws.Range("H79").Value = 725

$ws.Range("B83").Value = 19100
$ws.Range("C83").Value = 34
$ws.Range("D83").Value = 18228
$ws.Range("E83").Value = 752

$ws.Range("B98").Value = 9901
$ws.Range("C98").Value = 83
$ws.Range("D98").Value = 7358
$ws.Range("E98").Value = 2437
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 106

$ws.Range("B112").Value = 5701
$ws.Range("C112").Value = 4
$ws.Range("D112").Value = 3762
$ws.Range("E112").Value = 1761

$ws.Range("B119").Value = 4904
$ws.Range("C119").Value = 65
$ws.Range("D119").Value = 4294
$ws.Range("E119").Value = 564
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 46

$ws.Range("B123").Value = 4624
$ws.Range("C123").Value = 22
$ws.Range("D123").Value = 2767
$ws.Range("E123").Value = 1835
$ws.Range("H123").Value = 22

$ws.Range("B124").Value = 4611
$ws.Range("D124").Value = 3935
$ws.Range("E124").Value = 581
$ws.Range("H124").Value = 95

$ws.Range("B129").Value = 3569
$ws.Range("C129").Value = 130
$ws.Range("D129").Value = 1332
$ws.Range("E129").Value = 2098
$ws.Range("H129").Value = 139

$ws.Range("B130").Value = 3502
$ws.Range("C130").Value = 307
$ws.Range("D130").Value = 832
$ws.Range("E130").Value = 2635
$ws.Range("G130").Value = 3
$ws.Range("H130").Value = 35

$ws.Range("B131").Value = 3480
$ws.Range("C131").Value = 5
$ws.Range("D131").Value = 3315
$ws.Range("E131").Value = 107
$ws.Range("H131").Value = 58

$ws.Range("B138").Value = 3152
$ws.Range("C138").Value = 92
$ws.Range("D138").Value = 1610
$ws.Range("E138").Value = 1520
$ws.Range("G138").Value = 2
$ws.Range("H138").Value = 22

$ws.Range("B139").Value = 3080
$ws.Range("D139").Value = 837
$ws.Range("E139").Value = 2219
$ws.Range("H139").Value = 24

$ws.Range("B151").Value = 2126
$ws.Range("C151").Value = 15
$ws.Range("D151").Value = 1640
$ws.Range("E151").Value = 414

$ws.Range("B157").Value = 1595
$ws.Range("C157").Value = 17
$ws.Range("D157").Value = 1219
$ws.Range("E157").Value = 336

